# Update "想去人数" (want-to-go count) figures in column F across the
# four worksheets, matching the refreshed data snapshot from the
# gh-pages generator run at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value  = 258
$ws.Range("F5").Value  = 1028
$ws.Range("F6").Value  = 2265
$ws.Range("F10").Value = 197
$ws.Range("F11").Value = 154
$ws.Range("F13").Value = 56
$ws.Range("F14").Value = 93
$ws.Range("F15").Value = 1329
$ws.Range("F19").Value = 248

# --- Sheet: 演出 (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value  = 9
$ws.Range("F9").Value  = 119
$ws.Range("F11").Value = 35
$ws.Range("F12").Value = 210

# --- Sheet: 本地生活 (Local life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 1987
$ws.Range("F5").Value = 216

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value  = 1987
$ws.Range("F6").Value  = 216
$ws.Range("F12").Value = 258
$ws.Range("F13").Value = 1028
$ws.Range("F16").Value = 9
$ws.Range("F17").Value = 2265
$ws.Range("F18").Value = 119
$ws.Range("F21").Value = 35
$ws.Range("F24").Value = 197
$ws.Range("F25").Value = 210
$ws.Range("F26").Value = 154
$ws.Range("F28").Value = 56
$ws.Range("F29").Value = 93
$ws.Range("F31").Value = 1329
$ws.Range("F43").Value = 248
